$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 21: the Time cell (B21) had accidentally been entered as a date;
#     fix it to the actual time-range text. Style (s=2) stays as-is.
$ws.Range("B21").Value = "6:30 - 10:30"

# --- Row 23 (2/20 classtime entry): plain typed values into the existing
#     blank-template cells, so formatting stays at the template default
#     (s=4 / s=5) except the date cell, which gets date formatting below.
$ws.Range("B23").Value = "classtime"
$ws.Range("C23").Value = "classmates"
$ws.Range("D23").Value = "review last lecture contents, get feedback from last homework, learn new stuff"
$ws.Range("E23").Value = "learned more key expert practices, some ways of getting a higher level of abstractions of a system and some strategies of how code review works."
$ws.Range("F23").Value = "There are no simple plug-ins that can visualize an architecture pattern of a system, but there are other ways such as grouping source code and communications together based on folders."
$ws.Range("G23").Value = "feel terrified about having to document the architecture of Cassandra and other pieces of the homework"

# --- Row 22 (2/16 homework-3 entry): this one was produced by copying the
#     format of row 21 down onto row 22 first, then overwriting the values.
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B22").Value = "13:00 - 15:00"
$ws.Range("C22").Value = "self"
$ws.Range("D22").Value = "do homework 3"
$ws.Range("E22").Value = "Looked up stakeholders of Cassandra, essential functional and non-functional aspects of the system. "
$ws.Range("F22").Value = "It was not too hard to find stakeholders of Cassandra as it is a pretty well-known open source project. What interesting is that not just big companies are using it, even individual and some small businesses are using it as well.
"
$ws.Range("G22").Value = "feel accompolished after finishing most of the homework"

# --- Dates for the two new rows. Typing numbers directly keeps the
#     General style, so copy the date-formatted style from A21 first and
#     then write the serial value (mirrors the existing A21/A22/A23 pattern
#     instead of minting a brand-new date style).
$ws.Range("A21").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A22").Value = 43877
$ws.Range("A23").Value = 43881

# --- Row heights for the newly written rows (Excel auto-fit these once the
#     wrapped text was entered).
$ws.Rows.Item(22).RowHeight = 119.25
$ws.Rows.Item(23).RowHeight = 119
$ws.Rows.Item(24).RowHeight = 101

# --- Row 115 lost its "bottom block" shading/format (s=13/14) and picked up
#     the plain mid-table blank-row format (s=4/5) that row 114 uses.
$ws.Range("A114:G114").Copy()
$ws.Range("A115:G115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- A new blank row was appended at the very bottom of the table, matching
#     the formatting of the other trailing placeholder rows.
$ws.Range("A125:G125").Copy()
$ws.Range("A126:G126").PasteSpecial(-4122)
$excel.CutCopyMode = 0
